$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "wong3"

# Update cell values per row
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1

$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 6
$ws.Range("E3").Value = 14
$ws.Range("F3").Value = 14

$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 6

$ws.Range("B5").Value = 10
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 13

$ws.Range("E6").Value = 12
$ws.Range("F6").Value = 7

$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 16

$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 4
$ws.Range("E8").Value = 11
$ws.Range("F8").Value = 11

$ws.Range("E9").Value = 11
$ws.Range("F9").Value = 11

$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 4

$ws.Range("B13").Value = 6
$ws.Range("C13").Value = 6
$ws.Range("E13").Value = 14
$ws.Range("F13").Value = 14

$ws.Range("B14").Value = 10
$ws.Range("C14").Value = 10
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 4

$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 8

$ws.Range("E16").Value = 16
$ws.Range("F16").Value = 16

$ws.Range("B18").Value = 6
$ws.Range("C18").Value = 6
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 16

$ws.Range("B21").Value = 4
$ws.Range("C21").Value = 4
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = 8

$ws.Range("B22").Value = 7
$ws.Range("C22").Value = 7
$ws.Range("E22").Value = 28
$ws.Range("F22").Value = 26

$wb.Save()
